$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I10 was 0, user filled in an expense value of 20 ("Gas" column on 11/29)
$ws.Range("I10").Value = 20

# Extend the date column (A) down for the two new days, reusing the
# existing date-formatted style from A10 so no new number format is created.
$ws.Range("A10").Copy()
$ws.Range("A11:A12").PasteSpecial(-4122)
$ws.Range("A11").Value = 43799
$ws.Range("A12").Value = 43800

# Row 11 - 11/30/2019
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 16.5
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 26
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 2
$ws.Range("M11").Value = 3

# Row 12 - 12/1/2019
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 12.5
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 2
$ws.Range("M12").Value = 3

# Leave selection where the user ended up after entering the last row
$ws.Range("K12").Select()
